$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value corrections (missing-data imputation / reversal) ---

# Row 3 (RM 8): fill D3 with a recovered numeric reading
$ws.Range("D3").Value = -14.2

# Row 4 (RM 9): E4 becomes missing -> write as an empty text entry
# (a leading apostrophe forces a literal/text entry; an empty literal
# collapses to an empty string rather than clearing the cell outright)
$ws.Range("E4").Value = "'"
$ws.Range("E4").Style = "Normal"

# Row 5 (RM 14): D5 becomes missing
$ws.Range("D5").Value = "'"
$ws.Range("D5").Style = "Normal"

# Row 9 (RM 42): fill E9 with a recovered numeric reading
$ws.Range("E9").Value = -6.8

# Row 10 (RM 52 a): fill E10 with a recovered numeric reading
$ws.Range("E10").Value = -6.1

# Row 17 (RM 116): E17 becomes missing
$ws.Range("E17").Value = "'"
$ws.Range("E17").Style = "Normal"

# Row 18 (RM 120): E18 becomes missing
$ws.Range("E18").Value = "'"
$ws.Range("E18").Style = "Normal"

# Row 21 (RM 135): fill D21 with a recovered numeric reading
$ws.Range("D21").Value = -14.3

# Row 23 (RM 140): D23 becomes missing
$ws.Range("D23").Value = "'"
$ws.Range("D23").Style = "Normal"

# Row 34 (SC 193): fill D34 with a recovered numeric reading
$ws.Range("D34").Value = -14.7

# --- Remove rows: RM 232 (row 26) and SC 92 (row 28) ---
# Delete bottom-up so the remaining row indices don't shift underneath us.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
